$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new parameter row: LatticeScope_CH1_Mean with unit V
$ws.Range("A72").Value = "LatticeScope_CH1_Mean"
$ws.Range("B72").Value = "V"

# Widen column A slightly (onsite tweak while debugging the scope)
$ws.Columns.Item(1).ColumnWidth = 21.83

# Scroll down the sheet and move the selection near the newly added rows
[void]$ws.Range("C70").Select()
